$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# GRM_FOL_date: type changes from "date" to "text"
$ws.Range("C2").Value = "text"

# GRM_time_begin: type changes from "time" to "text"
$ws.Range("C4").Value = "text"

# GRM_time_end: type changes from "time" to "text"
$ws.Range("C6").Value = "text"
